$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "2020" column (N) mirroring the style of the existing
# "2019" column (M) for each data row, then fill in the 2020 values.

# Header row
$ws.Range("M4").Copy()
$ws.Range("N4").PasteSpecial(-4122)
$ws.Range("N4").Value = 2020

# Row 5
$ws.Range("M5").Copy()
$ws.Range("N5").PasteSpecial(-4122)
$ws.Range("N5").Value = 588.70000000000005

# Row 6
$ws.Range("M6").Copy()
$ws.Range("N6").PasteSpecial(-4122)
$ws.Range("N6").Value = 62.2

# Row 7
$ws.Range("M7").Copy()
$ws.Range("N7").PasteSpecial(-4122)
$ws.Range("N7").Value = 99.4

# Row 8
$ws.Range("M8").Copy()
$ws.Range("N8").PasteSpecial(-4122)
$ws.Range("N8").Value = 6.1

# Row 9 ("-" placeholder, no data)
$ws.Range("M9").Copy()
$ws.Range("N9").PasteSpecial(-4122)
$ws.Range("N9").Value = "-"

# Row 10 — value uses a new "0.0" number format instead of M10's style
$ws.Range("N10").Value = 71
$ws.Range("N10").NumberFormat = "0.0"

# Row 11
$ws.Range("M11").Copy()
$ws.Range("N11").PasteSpecial(-4122)
$ws.Range("N11").Value = 136.30000000000001

# Row 12
$ws.Range("M12").Copy()
$ws.Range("N12").PasteSpecial(-4122)
$ws.Range("N12").Value = 103.3

# Row 13
$ws.Range("M13").Copy()
$ws.Range("N13").PasteSpecial(-4122)
$ws.Range("N13").Value = 103.2

# Row 14
$ws.Range("M14").Copy()
$ws.Range("N14").PasteSpecial(-4122)
$ws.Range("N14").Value = 1.8

# Row 15 ("-" placeholder, no data)
$ws.Range("M15").Copy()
$ws.Range("N15").PasteSpecial(-4122)
$ws.Range("N15").Value = "-"

# Row 16
$ws.Range("M16").Copy()
$ws.Range("N16").PasteSpecial(-4122)
$ws.Range("N16").Value = 5.4

$excel.CutCopyMode = $false

# Move the active selection, matching the author's final cursor position.
$ws.Range("P15").Select() | Out-Null
